$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 2
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -8
$ws.Range("F12").Value = -7
